$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 47.67
$ws.Range("D5").Value = 2.29
$ws.Range("D6").Value = 3.29
$ws.Range("D10").Value = 3.59
$ws.Range("D11").Value = 15.88
